$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Germany Verbandsliga")

# --- Swap rows 45 and 46 (all columns B:AC; column A keeps the running id 43/44) ---
$r45 = $ws.Range("B45:AC45")
$r46 = $ws.Range("B46:AC46")
$v45 = $r45.Value()
$v46 = $r46.Value()
$r45.Value = $v46
$r46.Value = $v45

# --- Append new row 160 (new fixture) ---
# Copy formatting from the last existing row (159) so the new row matches the
# existing look & feel (bold/bordered id cell, date-formatted date cell, etc.)
$ws.Range("A159:AC159").Copy($ws.Range("A160:AC160"))

$ws.Cells.Item(160, 1).Value = 158
$ws.Cells.Item(160, 2).Value = 7873628
$ws.Cells.Item(160, 3).Value = "Germany Verbandsliga"
$ws.Cells.Item(160, 4).Value = "Germany Verbandsliga"
$ws.Cells.Item(160, 5).Value = 45346.5
$ws.Cells.Item(160, 6).Value = "SV Eintracht WaldMichelbach"
$ws.Cells.Item(160, 7).Value = "RotWeiss Frankfurt"
$ws.Cells.Item(160, 8).ClearContents()
$ws.Cells.Item(160, 9).ClearContents()
$ws.Cells.Item(160, 10).ClearContents()
$ws.Cells.Item(160, 11).Value = 1.909
$ws.Cells.Item(160, 12).Value = 3.75
$ws.Cells.Item(160, 13).Value = 3.1
$ws.Cells.Item(160, 14).Value = 1.909
$ws.Cells.Item(160, 15).Value = 3.75
$ws.Cells.Item(160, 16).Value = 3.2
$ws.Cells.Item(160, 17).Value = -0.5
$ws.Cells.Item(160, 18).Value = 1.95
$ws.Cells.Item(160, 19).Value = 1.85
$ws.Cells.Item(160, 20).Value = 3.25
$ws.Cells.Item(160, 21).Value = 1.85
$ws.Cells.Item(160, 22).Value = 1.95
$ws.Cells.Item(160, 23).Value = 0
$ws.Cells.Item(160, 24).Value = 0
$ws.Cells.Item(160, 25).Value = 0
$ws.Cells.Item(160, 26).Value = 0
$ws.Cells.Item(160, 27).Value = 0
$ws.Cells.Item(160, 28).ClearContents()
$ws.Cells.Item(160, 29).ClearContents()
